# ---------------------------------------------------------------------------
# Update "上海-漫展信息.xlsx" to match the scraped output generated at 456a3b4.
#
# Changes:
#  1. Sheet "展览"  (Exhibitions)  - "want to go" counts (column F) bumped.
#  2. Sheet "本地生活" (Local life) - "want to go" counts (column F) bumped.
#  3. Sheet "演出"  (Performances) - "want to go" counts bumped AND a new
#     event row is inserted at row 37 ("【浪漫520】遇见爱情主题视听音乐会
#     《一生所爱》"), pushing the previous rows 37-47 down to 38-48.
#  4. Sheet "全部类型" (All types, a combined listing) - "want to go" counts
#     bumped (no row insert needed there).
#
# NOTE: this COM runtime does not correctly bind COM objects passed through
# *named* function parameters (eg. "Test-Func -ws $ws") -- they arrive as
# $null inside the function body. Positional parameter passing works fine,
# so all helper-function calls below use positional arguments only.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param($ws, $map)
    foreach ($row in $map.Keys) {
        $ws.Range("F$row").Value = $map[$row]
    }
}

# ---------------------------------------------------------------------------
# 1. Sheet "展览"
# ---------------------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$expoFChanges = @{
    2 = 257; 3 = 881; 4 = 553; 5 = 2297; 7 = 122; 9 = 1172; 10 = 1047; 11 = 3058;
    14 = 1105; 16 = 534; 17 = 240; 18 = 615; 19 = 1131; 20 = 1131; 21 = 163;
    23 = 194; 24 = 539; 25 = 252; 26 = 648; 28 = 10; 29 = 843; 30 = 89; 32 = 55;
    33 = 1057; 34 = 5071; 35 = 529; 36 = 263; 37 = 139; 39 = 7
}
Set-FValues $wsExpo $expoFChanges

# ---------------------------------------------------------------------------
# 2. Sheet "本地生活"
# ---------------------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$localFChanges = @{ 5 = 439; 6 = 423 }
Set-FValues $wsLocal $localFChanges

# ---------------------------------------------------------------------------
# 3. Sheet "演出"
# ---------------------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

# 3a. "want to go" bumps for rows that are NOT affected by the later insert
#     (these rows sit above row 37, so their row numbers do not change).
$showFChangesPreInsert = @{ 6 = 414; 21 = 46; 25 = 389; 30 = 1; 34 = 63 }
Set-FValues $wsShow $showFChangesPreInsert

# 3b. Insert the new event as row 37; everything from old row 37 on down
#     shifts to row+1.
$wsShow.Rows(37).Insert()

# Copy formatting (bold/centered/bordered index-number style) from the cell
# that used to be A37 (now sitting at A38) onto the freshly inserted A37.
$wsShow.Range("A38").Copy()
$wsShow.Range("A37").PasteSpecial(-4122)   # xlPasteFormats
$wsShow.Application.CutCopyMode = $false

# New row content.
$wsShow.Range("A37").Value = 36
$wsShow.Range("B37").Value = "2024-05-19"
$wsShow.Range("C37").Value = "上海·【浪漫520】遇见爱情主题视听音乐会《一生所爱》"
$wsShow.Range("D37").Value = "岳阳街道人民南路 69号 云间剧院"
$wsShow.Range("E37").Value = "2024.05.19 19:30-05.19 21:00"
$wsShow.Range("F37").Value = 0
$wsShow.Range("G37").Value = 60
$wsShow.Range("H37").Value = "https://show.bilibili.com/platform/detail.html?id=84083"
$wsShow.Range("I37").Value = "//i0.hdslb.com/bfs/openplatform/202404/p9zrWxJ11712805726433.png"

# 3c. The "A" (index) column is always row-number minus one; after the row
#     insert, Excel shifted the literal numbers down along with everything
#     else, so they need to be re-stamped for every row from the inserted
#     one through the end of the sheet.
$lastRow = $wsShow.Range("A1").CurrentRegion.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $wsShow.Range("A$r").Value = $r - 1
}

# 3d. "want to go" bumps that land on the shifted rows (38 and 40 get an
#     extra bump on top of the plain shift down from 37 and 39).
$showFChangesPostInsert = @{ 38 = 446; 40 = 14 }
Set-FValues $wsShow $showFChangesPostInsert

# ---------------------------------------------------------------------------
# 4. Sheet "全部类型"
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$allFChanges = @{
    3 = 257; 4 = 439; 6 = 881; 7 = 414; 8 = 553; 9 = 2297; 11 = 122; 13 = 1172;
    15 = 1047; 16 = 3058; 20 = 1105; 21 = 423; 22 = 534; 23 = 240; 24 = 615;
    25 = 1131; 26 = 1131; 27 = 163; 28 = 46; 31 = 194; 32 = 252; 34 = 648;
    36 = 389; 38 = 843; 39 = 89; 42 = 55; 43 = 1058; 44 = 5071; 45 = 63;
    46 = 529; 47 = 446; 48 = 446; 49 = 263
}
Set-FValues $wsAll $allFChanges

$wb.Save()
